# Generate Report for Handback
# - Marks the two files as "Handed back: in sync with en-US"
# - Records the latest handback target/file/datetime for the zh-cn and de-de
#   variant sheets, including hyperlinks to the handed-back target files
# - Widens a few columns on the Overview / zh-cn / de-de sheets so the new
#   content is readable

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$zhDateTime = "2016-08-17 22:23:51"
$deDateTime = "2016-08-17 22:23:58"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1dd0d453c256f0f22309e4370af57139bd19b5c/e2e/"
$file1Md = "121b94b7-67fd-4439-a95e-b403e815b5ae.md"
$file2Md = "ac3fb723-3934-4251-86ab-b855325dbdf1.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("K2").Value = $zhDateTime
$wsZh.Range("K3").Value = $zhDateTime

$wsZh.Range("J2").Value = "121b94b7-67fd-4439-a95e-b403e815b5ae.fd59b22916b813b99d79b799f24e4722797f33c2.zh-cn.xlf"
$wsZh.Range("J3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.4ea0bbd83bb899adbf7afd34e2211aa107e5fa86.zh-cn.xlf"

$wsZh.Range("I2").Value = $file1Md
$wsZh.Range("I3").Value = $file2Md

$zhLinks = @($wsZh.Hyperlinks)
$zhLinks[1].Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($repoBase + $file1Md), "", "", $file1Md)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($repoBase + $file2Md), "", "", $file2Md)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($repoBase + $file2Md), "", "", $file2Md)

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("K2").Value = $deDateTime
$wsDe.Range("K3").Value = $deDateTime

$wsDe.Range("J2").Value = "121b94b7-67fd-4439-a95e-b403e815b5ae.fd59b22916b813b99d79b799f24e4722797f33c2.de-de.xlf"
$wsDe.Range("J3").Value = "ac3fb723-3934-4251-86ab-b855325dbdf1.4ea0bbd83bb899adbf7afd34e2211aa107e5fa86.de-de.xlf"

$wsDe.Range("I2").Value = $file1Md
$wsDe.Range("I3").Value = $file2Md

$deLinks = @($wsDe.Hyperlinks)
$deLinks[1].Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($repoBase + $file1Md), "", "", $file1Md)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($repoBase + $file2Md), "", "", $file2Md)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($repoBase + $file2Md), "", "", $file2Md)

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# Overview sheet - the zh-cn / de-de summary columns get wider too
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527
